{"js": "// Replace the date line and each two-digit multiplication expression in\n// the worksheet with the updated values from the new day's output.\nconst replacements = [\n  [\"2025-02-03 Monday\", \"2025-02-04 Tuesday\"],\n  [\"99\u00d721=\", \"99\u00d720=\"],\n  [\"88\u00d791=\", \"24\u00d790=\"],\n  [\"33\u00d732=\", \"33\u00d737=\"],\n  [\"12\u00d769=\", \"33\u00d798=\"],\n  [\"91\u00d787=\", \"93\u00d778=\"],\n  [\"71\u00d789=\", \"75\u00d722=\"],\n  [\"88\u00d741=\", \"55\u00d786=\"],\n  [\"82\u00d782=\", \"27\u00d778=\"],\n  [\"20\u00d723=\", \"56\u00d787=\"],\n  [\"43\u00d714=\", \"64\u00d792=\"],\n  [\"28\u00d781=\", \"49\u00d740=\"],\n  [\"49\u00d760=\", \"57\u00d716=\"],\n  [\"32\u00d782=\", \"97\u00d787=\"],\n  [\"67\u00d753=\", \"46\u00d757=\"],\n  [\"66\u00d744=\", \"42\u00d719=\"],\n  [\"17\u00d788=\", \"70\u00d725=\"],\n  [\"86\u00d725=\", \"28\u00d786=\"],\n  [\"97\u00d772=\", \"95\u00d762=\"],\n  [\"28\u00d756=\", \"64\u00d762=\"],\n  [\"55\u00d773=\", \"83\u00d728=\"],\n  [\"62\u00d717=\", \"24\u00d741=\"],\n  [\"31\u00d747=\", \"92\u00d760=\"],\n  [\"14\u00d719=\", \"75\u00d716=\"],\n  [\"60\u00d747=\", \"90\u00d763=\"],\n  [\"95\u00d754=\", \"20\u00d789=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit multiplication expression in\n# the worksheet with the updated values from the new day's output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-03 Monday\", \"2025-02-04 Tuesday\"),\n    @(\"99\u00d721=\", \"99\u00d720=\"),\n    @(\"88\u00d791=\", \"24\u00d790=\"),\n    @(\"33\u00d732=\", \"33\u00d737=\"),\n    @(\"12\u00d769=\", \"33\u00d798=\"),\n    @(\"91\u00d787=\", \"93\u00d778=\"),\n    @(\"71\u00d789=\", \"75\u00d722=\"),\n    @(\"88\u00d741=\", \"55\u00d786=\"),\n    @(\"82\u00d782=\", \"27\u00d778=\"),\n    @(\"20\u00d723=\", \"56\u00d787=\"),\n    @(\"43\u00d714=\", \"64\u00d792=\"),\n    @(\"28\u00d781=\", \"49\u00d740=\"),\n    @(\"49\u00d760=\", \"57\u00d716=\"),\n    @(\"32\u00d782=\", \"97\u00d787=\"),\n    @(\"67\u00d753=\", \"46\u00d757=\"),\n    @(\"66\u00d744=\", \"42\u00d719=\"),\n    @(\"17\u00d788=\", \"70\u00d725=\"),\n    @(\"86\u00d725=\", \"28\u00d786=\"),\n    @(\"97\u00d772=\", \"95\u00d762=\"),\n    @(\"28\u00d756=\", \"64\u00d762=\"),\n    @(\"55\u00d773=\", \"83\u00d728=\"),\n    @(\"62\u00d717=\", \"24\u00d741=\"),\n    @(\"31\u00d747=\", \"92\u00d760=\"),\n    @(\"14\u00d719=\", \"75\u00d716=\"),\n    @(\"60\u00d747=\", \"90\u00d763=\"),\n    @(\"95\u00d754=\", \"20\u00d789=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
